$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) - force text format to avoid Excel auto-numeric conversion,
# then clear formatting so the cell style matches the original (no explicit style index).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.367.02"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.845.75"
$ws.Range("D3").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.63"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6275"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("D7").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07487"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2892"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.37"
$ws.Range("D10").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.846.32"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.002"
$ws.Range("D13").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6788"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001028"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.65"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.110.10"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.092"
$ws.Range("D18").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.392.14"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.15"
$ws.Range("D20").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.418"
$ws.Range("D23").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.01"
$ws.Range("D25").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.389"
$ws.Range("D27").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.109"
$ws.Range("D32").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6929"
$ws.Range("D36").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.588"
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.852"
$ws.Range("D38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.251.25"
$ws.Range("D39").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.507"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9039"
$ws.Range("D42").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.010.43"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.28"
$ws.Range("D45").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.078"
$ws.Range("D47").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.993"
$ws.Range("D49").ClearFormats()

# Update Volume(1h) column (E) - plain percentage text with padding spaces.
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -2.72%  "
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("E10").Value = "  -2.53%  "
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("E15").Value = "  -5.44%  "
$ws.Range("E16").Value = "  -1.13%  "
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("E21").Value = "  -0.86%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("E28").Value = "  -0.83%  "
$ws.Range("E29").Value = "  +3.41%  "
$ws.Range("E30").Value = "  +0.79%  "
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("E34").Value = "  -1.39%  "
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("E38").Value = "  +3.61%  "
$ws.Range("E40").Value = "  +1.48%  "
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("E47").Value = "  -1.50%  "
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("E51").Value = "  -4.85%  "

Write-Host "Crypto price/volume updates applied"
